# Update the error-metric table (ME, MAE, MSE, RMSE, SE, N) in the "first"
# sheet with recalculated values that reflect one additional ifo GDP
# component observation being included in the preprocessing, per the
# commit "added ifo gdp component analysis preprocessing".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1279039517177195, 1.47408018424297,   8.675642101351979, 2.945444296087091, 3.008801495394956, 23),
    @(0.04717506928171853,1.491543766192013,  8.455526954849582, 2.907838880483164, 2.975876279823447, 22),
    @(-0.4214464237818173,1.016004021736437,  3.540968988508833, 1.881746260394539, 1.879233752200579, 21),
    @(-0.04213737398666674,0.7820104901831507,1.001266313089072, 1.000632956227743, 1.025717081543998, 20),
    @(-0.0547463289693642,0.7907933146236926, 1.075494849714573, 1.037060677932865, 1.063992901354768, 19),
    @(-0.1059385590168174, 0.6529650943868757,0.7643346331479901,0.8742623365718039,0.8929794736955126, 18),
    @(-0.06115313360779569,0.7409233030117504,0.8853216378027108,0.9409153191455174,0.9678227124871316, 17),
    @(0.04633068755579241, 0.7066816895393113,0.7257841232605365,0.8519296468961135,0.8785670701173324, 16),
    @(-0.03682260890791801,0.5844290800941564,0.7081477771230911,0.841515167494378, 0.8702166433441642, 15),
    @(-0.00294460816098668,0.5209809220943237,0.4617339020710022,0.6795100455997706,0.7051542787521982, 14)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
